# Update the cryptocurrency price/volume table to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" column (D) values, keyed by row number.
$priceUpdates = @{
    2 = '70.677.77'
    3 = '3.518.52'
    4 = '0.999'
    5 = '618.63'
    6 = '173.00'
    8 = '3.512.11'
    12 = '0.583'
    13 = '46.37'
    15 = '4.091.41'
    17 = '608.05'
    18 = '3.523.37'
    19 = '70.809.58'
    20 = '0.121'
    21 = '17.70'
    22 = '0.881'
    23 = '9.19'
    24 = '97.98'
    25 = '15.61'
    29 = '33.65'
    30 = '9.04'
    31 = '3.00'
    32 = '8.08'
    34 = '640.14'
    36 = '0.0997'
    37 = '10.79'
    40 = '1.00'
    41 = '56.58'
    43 = '3.351.63'
    44 = '0.0₃0716'
    45 = '0.311'
    46 = '2.92'
    47 = '31.75'
    48 = '2.53'
    50 = '134.70'
}

# New "Volume(1h)" column (E) values, keyed by row number.
$volumeUpdates = @{
    2 = '  -0.22%  '
    3 = '  -2.21%  '
    4 = '  -0.03%  '
    5 = '  +2.27%  '
    6 = '  -0.81%  '
    7 = '  -0.96%  '
    8 = '  -2.19%  '
    9 = '  +0.02%  '
    10 = '  -2.11%  '
    11 = '  -5.02%  '
    12 = '  -1.41%  '
    13 = '  -1.82%  '
    14 = '  -1.53%  '
    15 = '  -2.00%  '
    16 = '  -1.34%  '
    17 = '  -1.72%  '
    18 = '  -1.98%  '
    19 = '  -0.17%  '
    20 = '  +1.01%  '
    21 = '  +0.87%  '
    22 = '  -1.33%  '
    23 = '  -0.50%  '
    24 = '  -0.08%  '
    25 = '  -3.34%  '
    26 = '  -2.47%  '
    27 = '  +0.12%  '
    28 = '  -4.05%  '
    29 = '  -2.12%  '
    30 = '  -3.17%  '
    31 = '  -3.15%  '
    32 = '  -5.60%  '
    33 = '  -1.23%  '
    34 = '  +2.60%  '
    35 = '  -6.82%  '
    36 = '  -2.70%  '
    37 = '  -1.07%  '
    38 = '  -0.45%  '
    39 = '  -9.25%  '
    40 = '  +0.27%  '
    41 = '  -1.91%  '
    42 = '  -1.71%  '
    43 = '  -1.57%  '
    44 = '  -0.56%  '
    45 = '  -4.57%  '
    46 = '  -3.28%  '
    47 = '  -4.10%  '
    48 = '  -6.81%  '
    49 = '  -0.64%  '
    50 = '  +1.45%  '
}

# Rows whose ranking order changed: new Coin name (B) and Link (C) values.
$coinUpdates = @{
    24 = 'Litecoin'
    25 = 'InternetComputer(DFINITY)'
    40 = 'FirstDigitalUSD'
    41 = 'OKB'
}
$linkUpdates = @{
    24 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    25 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    40 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    41 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
}

# Rows whose Price text would otherwise be auto-converted by Excel into a
# number/date (plain decimals, or the PEPE row's subscript-digit notation).
# A leading apostrophe forces Excel to keep the literal text, exactly as if
# a person had typed it into the cell.
$quotePrefixRows = @(4, 5, 6, 12, 13, 17, 20, 21, 22, 23, 24, 25, 29, 30, 31, 32, 34, 36, 37, 40, 41, 44, 45, 46, 47, 48, 50)

foreach ($row in $priceUpdates.Keys) {
    $value = $priceUpdates[$row]
    if ($quotePrefixRows -contains $row) {
        $value = "'" + $value
    }
    $ws.Range("D$row").Value = $value
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

foreach ($row in $coinUpdates.Keys) {
    $ws.Range("B$row").Value = $coinUpdates[$row]
}

foreach ($row in $linkUpdates.Keys) {
    $ws.Range("C$row").Value = $linkUpdates[$row]
}
